$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 9.805277777777777

for ($r = 31; $r -le 53; $r++) {
    $ws.Range("I$r").Value = $newValue
}
